$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto Price (D) / Volume(1h) (E) columns to the latest scrape,
# per commit "Updated cryptos list on Mon Dec 25 16:34:44 UTC 2023 with GitHub Actions".
# Both columns hold plain text in the source workbook (note values such as
# "43.802.81" using "." as a thousands separator, or "0.640"/"42.50" with
# significant trailing zeros). A leading apostrophe forces Excel to keep
# number-like strings as literal text instead of silently re-parsing /
# renormalising them (which would drop the trailing zeros or coerce to a number).

$ws.Range("D2").Value = "43.802.81"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.282.69"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'123.46"
$ws.Range("E5").Value = "  +7.97%  "
$ws.Range("D6").Value = "'266.21"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").Value = "'0.640"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.627"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "'48.52"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").Value = "'0.0944"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'9.06"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'15.53"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "'0.902"
$ws.Range("E15").Value = "  +4.95%  "
$ws.Range("D16").Value = "2.627.10"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "2.277.52"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "43.744.96"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'7.03"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "'72.49"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'2.45"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'236.15"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").Value = "'9.58"
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").Value = "'11.87"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "'42.50"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "'3.37"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'172.18"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").Value = "'21.73"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").Value = "'5.76"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").Value = "'0.0380"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").Value = "'4.71"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "'4.11"
$ws.Range("E38").Value = "  +7.58%  "
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  +6.19%  "
$ws.Range("D41").Value = "'75.81"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("D42").Value = "'13.92"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "'0.239"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'5.85"
$ws.Range("E45").Value = "  -8.36%  "
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "'75.87"
$ws.Range("E47").Value = "  +41.59%  "
$ws.Range("D48").Value = "'1.27"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "'8.60"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").Value = "'102.05"
$ws.Range("E51").Value = "  -0.61%  "
